$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the two new SDL error code definitions (rows 83 and 84).
# Columns A, C, D, E, F, G, H are driven by existing shared formulas
# (si="17".."23") and recompute automatically once B/I are populated.
$ws.Range("B83").Value = "SDL_CREATED"
$ws.Range("I83").Value = "S"

$ws.Range("B84").Value = "SDL_NOTCREATED"
$ws.Range("I84").Value = "S"

# Move the active selection to the newly filled rows, matching the
# author's saved cursor position (A83, selection A83:C84).
$excel.Goto($ws.Range("A83:C84"))
